# Fruta / hortaliza, semanal
# The underlying source data for this market/category ("Achicoria") was
# refreshed; rows 2-24 were re-pulled and now line up with different
# dates (column D) and therefore different Volumen/Precio figures
# (columns J, K, L, M, P). Row 25 is unchanged.
#
# Capture the "before" values for the columns that vary per-row, then
# write the "after" (new) values back into the same row positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 24

# Snapshot the original (pre-edit) values for the columns that change.
$orig = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $orig[$r] = @{
        D = $ws.Range("D$r").Value2
        J = $ws.Range("J$r").Value2
        K = $ws.Range("K$r").Value2
        L = $ws.Range("L$r").Value2
        M = $ws.Range("M$r").Value2
        P = $ws.Range("P$r").Value2
    }
}

# Maps each destination row to the source row whose original data it
# should now display.
$rowMap = @{
    2  = 9
    3  = 8
    4  = 18
    5  = 17
    6  = 23
    7  = 13
    8  = 5
    9  = 14
    10 = 7
    11 = 16
    12 = 21
    13 = 3
    14 = 22
    15 = 20
    16 = 11
    17 = 24
    18 = 19
    19 = 15
    20 = 4
    21 = 2
    22 = 10
    23 = 12
    24 = 6
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $vals = $orig[$srcRow]

    $ws.Range("D$destRow").Value = $vals.D
    $ws.Range("J$destRow").Value = $vals.J
    $ws.Range("K$destRow").Value = $vals.K
    $ws.Range("L$destRow").Value = $vals.L
    $ws.Range("M$destRow").Value = $vals.M
    $ws.Range("P$destRow").Value = $vals.P
}
